$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 16) "320018701475"
Set-TextValue $ws.Cells.Item(2, 17) "`$20.36"
$ws.Cells.Item(2, 18).Value = "FAIL"

# Row 3
Set-TextValue $ws.Cells.Item(3, 16) "320018701497"
Set-TextValue $ws.Cells.Item(3, 17) "`$29.35"
$ws.Cells.Item(3, 18).Value = "FAIL"

# Row 4
Set-TextValue $ws.Cells.Item(4, 16) "320018701523"
Set-TextValue $ws.Cells.Item(4, 17) "`$33.84"
$ws.Cells.Item(4, 18).Value = "FAIL"

# Row 5
Set-TextValue $ws.Cells.Item(5, 16) "320018701545"
Set-TextValue $ws.Cells.Item(5, 17) "`$46.27"
$ws.Cells.Item(5, 18).Value = "FAIL"

# Row 6
Set-TextValue $ws.Cells.Item(6, 16) "320018701589"
Set-TextValue $ws.Cells.Item(6, 17) "`$60.01"
$ws.Cells.Item(6, 18).Value = "FAIL"

# Row 7
Set-TextValue $ws.Cells.Item(7, 16) "320018701604"
Set-TextValue $ws.Cells.Item(7, 17) "`$236.37"
$ws.Cells.Item(7, 18).Value = "FAIL"

# Row 8
Set-TextValue $ws.Cells.Item(8, 16) "320018701648"
Set-TextValue $ws.Cells.Item(8, 17) "`$20.36"
$ws.Cells.Item(8, 18).Value = "FAIL"

# Row 9
Set-TextValue $ws.Cells.Item(9, 16) "320018702081"
Set-TextValue $ws.Cells.Item(9, 17) "`$24.85"
$ws.Cells.Item(9, 18).Value = "FAIL"

# Row 10
Set-TextValue $ws.Cells.Item(10, 16) "320018702130"
Set-TextValue $ws.Cells.Item(10, 17) "`$29.35"
$ws.Cells.Item(10, 18).Value = "FAIL"

# Row 11
Set-TextValue $ws.Cells.Item(11, 16) "320018702162"
Set-TextValue $ws.Cells.Item(11, 17) "`$43.09"
$ws.Cells.Item(11, 18).Value = "FAIL"

# Row 12
Set-TextValue $ws.Cells.Item(12, 16) "320018702210"
Set-TextValue $ws.Cells.Item(12, 17) "`$56.58"
$ws.Cells.Item(12, 18).Value = "FAIL"

# Row 13
Set-TextValue $ws.Cells.Item(13, 16) "320018702232"
Set-TextValue $ws.Cells.Item(13, 17) "`$15.86"
$ws.Cells.Item(13, 18).Value = "FAIL"

# Row 14
Set-TextValue $ws.Cells.Item(14, 16) "320018702265"
Set-TextValue $ws.Cells.Item(14, 17) "`$19.30"
$ws.Cells.Item(14, 18).Value = "FAIL"

# Row 15
Set-TextValue $ws.Cells.Item(15, 16) "320018702287"
Set-TextValue $ws.Cells.Item(15, 17) "`$22.74"
$ws.Cells.Item(15, 18).Value = "FAIL"

# Row 16
Set-TextValue $ws.Cells.Item(16, 16) "320018702324"
Set-TextValue $ws.Cells.Item(16, 17) "`$33.84"
$ws.Cells.Item(16, 18).Value = "FAIL"

# Row 17
Set-TextValue $ws.Cells.Item(17, 16) "320018702368"
Set-TextValue $ws.Cells.Item(17, 17) "`$45.21"
$ws.Cells.Item(17, 18).Value = "FAIL"

# Row 18
Set-TextValue $ws.Cells.Item(18, 16) "320018702405"
Set-TextValue $ws.Cells.Item(18, 17) "`$46.27"
$ws.Cells.Item(18, 18).Value = "FAIL"

# Row 19
Set-TextValue $ws.Cells.Item(19, 16) "320018702427"
Set-TextValue $ws.Cells.Item(19, 17) "`$57.63"
$ws.Cells.Item(19, 18).Value = "FAIL"

# Row 20
Set-TextValue $ws.Cells.Item(20, 16) "320018702450"
Set-TextValue $ws.Cells.Item(20, 17) "`$66.89"
$ws.Cells.Item(20, 18).Value = "FAIL"

# Row 21
Set-TextValue $ws.Cells.Item(21, 16) "320018702471"
Set-TextValue $ws.Cells.Item(21, 17) "`$118.70"
$ws.Cells.Item(21, 18).Value = "FAIL"

# Row 22
Set-TextValue $ws.Cells.Item(22, 16) "320018702508"

# Row 23
Set-TextValue $ws.Cells.Item(23, 16) "320018702519"

# Row 24
Set-TextValue $ws.Cells.Item(24, 16) "320018702541"

# Row 25
Set-TextValue $ws.Cells.Item(25, 16) "320018702563"

# Row 26
Set-TextValue $ws.Cells.Item(26, 16) "320018702574"
